$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.735.55"
$ws.Range("E2").Value = "  -1.50%  "
$ws.Range("D3").Value = "3.482.80"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'581.95"
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("D6").Value = "'129.90"
$ws.Range("E6").Value = "  -3.02%  "
$ws.Range("D7").Value = "3.481.85"
$ws.Range("E7").Value = "  -0.88%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -2.10%  "
$ws.Range("D10").Value = "'0.123"
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("D11").Value = "'7.09"
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").Value = "'0.376"
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("D13").Value = "4.055.17"
$ws.Range("E13").Value = "  -1.51%  "
$ws.Range("D14").Value = "'27.15"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").Value = "3.484.15"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("E17").Value = "  -3.44%  "
$ws.Range("D18").Value = "63.696.67"
$ws.Range("E18").Value = "  -1.66%  "
$ws.Range("D19").Value = "'9.82"
$ws.Range("E19").Value = "  -1.91%  "
$ws.Range("D20").Value = "'14.02"
$ws.Range("E20").Value = "  -2.46%  "
$ws.Range("D21").Value = "'5.59"
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("D22").Value = "'379.34"
$ws.Range("E22").Value = "  -3.31%  "
$ws.Range("D23").Value = "'0.571"
$ws.Range("E23").Value = "  -1.22%  "
$ws.Range("D24").Value = "3.619.35"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").Value = "'73.06"
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("D28").Value = "'1.56"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "'7.44"
$ws.Range("E30").Value = "  -2.61%  "
$ws.Range("D31").Value = "'8.18"
$ws.Range("E31").Value = "  -1.37%  "
$ws.Range("D32").Value = "'2.21"
$ws.Range("E32").Value = "  -3.16%  "
$ws.Range("D33").Value = "3.488.66"
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("D35").Value = "'23.30"
$ws.Range("E35").Value = "  -3.58%  "
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("D37").Value = "'5.25"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "'6.89"
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("E39").Value = "  -1.71%  "
$ws.Range("D40").Value = "'159.89"
$ws.Range("E40").Value = "  -5.24%  "
$ws.Range("D41").Value = "'0.0788"
$ws.Range("E41").Value = "  -3.47%  "
$ws.Range("D42").Value = "'0.808"
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("D43").Value = "'25.90"
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("D44").Value = "'0.998"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").Value = "'41.68"
$ws.Range("E45").Value = "  -2.33%  "
$ws.Range("D46").Value = "'1.20"
$ws.Range("E46").Value = "  -3.94%  "
$ws.Range("D47").Value = "'4.34"
$ws.Range("E47").Value = "  -1.70%  "
$ws.Range("D48").Value = "'1.60"
$ws.Range("E48").Value = "  -2.45%  "
$ws.Range("D49").Value = "2.417.07"
$ws.Range("E49").Value = "  +1.44%  "
$ws.Range("D50").Value = "'6.79"
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("D51").Value = "'0.884"
$ws.Range("E51").Value = "  -1.13%  "
